$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row + data cells (D/E/F column realignment, case normalisation) ---
$ws.Cells.Item(1,1).Value = "industry"
$ws.Cells.Item(1,2).Value = "unit"
$ws.Cells.Item(1,3).Value = "process"
$ws.Cells.Item(1,4).Value = "carbon (kg CO2 eq)"
$ws.Cells.Item(1,5).Value = "ced (MJ)"
$ws.Cells.Item(1,6).Value = "climate change (kg CO2 eq)"
$ws.Cells.Item(1,7).Value = "region"
$ws.Cells.Item(2,4).Value = 2858.327333333334
$ws.Cells.Item(2,5).Value = 46758.303
$ws.Cells.Item(2,6).Value = 0.079697751
$ws.Cells.Item(3,4).Value = 2799.252866666667
$ws.Cells.Item(3,5).Value = 46790.267
$ws.Cells.Item(3,6).Value = 0.078050598
$ws.Cells.Item(4,4).Value = 748.278
$ws.Cells.Item(4,5).Value = 47063.868
$ws.Cells.Item(4,6).Value = 0.020863977
$ws.Cells.Item(5,4).Value = 294.6922866666667
$ws.Cells.Item(5,5).Value = 38740.507
$ws.Cells.Item(5,6).Value = 0.0082168028
$ws.Cells.Item(6,4).Value = 2880.933333333333
$ws.Cells.Item(6,5).Value = 44226.866
$ws.Cells.Item(6,6).Value = 0.080328065
$ws.Cells.Item(7,4).Value = 2790.875466666667
$ws.Cells.Item(7,5).Value = 43989.878
$ws.Cells.Item(7,6).Value = 0.077817014
$ws.Cells.Item(8,4).Value = 1436.143666666667
$ws.Cells.Item(8,5).Value = 41236.976
$ws.Cells.Item(8,6).Value = 0.040043496
$ws.Cells.Item(9,4).Value = 3.971818133333334
$ws.Cells.Item(9,5).Value = 58.84974
$ws.Cells.Item(9,6).Value = 0.00011074483
$ws.Cells.Item(10,4).Value = 3.971818133333334
$ws.Cells.Item(10,5).Value = 58.84974
$ws.Cells.Item(10,6).Value = 0.00011074483
$ws.Cells.Item(11,4).Value = 1714.285733333333
$ws.Cells.Item(11,5).Value = 26916.372
$ws.Cells.Item(11,6).Value = 0.047798834
$ws.Cells.Item(12,4).Value = 11.76893066666667
$ws.Cells.Item(12,5).Value = 174.91564
$ws.Cells.Item(12,6).Value = 0.00032814902
$ws.Cells.Item(13,4).Value = 0.1367532266666667
$ws.Cells.Item(13,5).Value = 2.1555017
$ws.Cells.Item(13,6).Value = 0.0000038130428
$ws.Cells.Item(14,4).Value = 0.1479167533333333
$ws.Cells.Item(14,5).Value = 2.331461
$ws.Cells.Item(14,6).Value = 0.0000041243116
$ws.Cells.Item(15,4).Value = 0.15908028
$ws.Cells.Item(15,5).Value = 2.5074203
$ws.Cells.Item(15,6).Value = 0.0000044355804
$ws.Cells.Item(16,4).Value = 0.1869891066666667
$ws.Cells.Item(16,5).Value = 2.9473186
$ws.Cells.Item(16,6).Value = 0.0000052137524
$ws.Cells.Item(17,4).Value = 0.2455976266666667
$ws.Cells.Item(17,5).Value = 3.871105
$ws.Cells.Item(17,6).Value = 0.0000068479136
$ws.Cells.Item(18,4).Value = 0.1233076533333333
$ws.Cells.Item(18,5).Value = 2.0577513
$ws.Cells.Item(18,6).Value = 0.0000034381445
$ws.Cells.Item(19,4).Value = 0.1317149933333333
$ws.Cells.Item(19,5).Value = 2.1980525
$ws.Cells.Item(19,6).Value = 0.0000036725634
$ws.Cells.Item(20,4).Value = 0.1401223333333333
$ws.Cells.Item(20,5).Value = 2.3383537
$ws.Cells.Item(20,6).Value = 0.0000039069824
$ws.Cells.Item(21,4).Value = 0.1989737066666667
$ws.Cells.Item(21,5).Value = 3.3204623
$ws.Cells.Item(21,6).Value = 0.000005547915
$ws.Cells.Item(22,4).Value = 0.04554260866666666
$ws.Cells.Item(22,5).Value = 1.196605
$ws.Cells.Item(22,6).Value = 0.0000012698488
$ws.Cells.Item(23,4).Value = 0.051072782
$ws.Cells.Item(23,5).Value = 1.3419071
$ws.Cells.Item(23,6).Value = 0.0000014240447
$ws.Cells.Item(24,4).Value = 0.05513908666666667
$ws.Cells.Item(24,5).Value = 1.4487468
$ws.Cells.Item(24,6).Value = 0.000001537424
$ws.Cells.Item(25,4).Value = 0.05465113
$ws.Cells.Item(25,5).Value = 1.4359261
$ws.Cells.Item(25,6).Value = 0.0000015238185
$ws.Cells.Item(26,4).Value = 0.07221756666666668
$ws.Cells.Item(26,5).Value = 1.8974737
$ws.Cells.Item(26,6).Value = 0.0000020136173
$ws.Cells.Item(27,4).Value = 0
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(29,4).Value = 1.304549666666667
$ws.Cells.Item(29,5).Value = 51.774007
$ws.Cells.Item(29,6).Value = 0.000036374306
$ws.Cells.Item(30,4).Value = 0.8929654666666667
$ws.Cells.Item(30,5).Value = 50.491708
$ws.Cells.Item(30,6).Value = 0.000024898247
$ws.Cells.Item(31,4).Value = 3.679221
$ws.Cells.Item(31,5).Value = 53.780909
$ws.Cells.Item(31,6).Value = 0.00010258644
$ws.Cells.Item(32,4).Value = 0.8410644666666667
$ws.Cells.Item(32,5).Value = 52.899788
$ws.Cells.Item(32,6).Value = 0.00002345111
$ws.Cells.Item(33,4).Value = 1.4478264
$ws.Cells.Item(33,5).Value = 51.161439
$ws.Cells.Item(33,6).Value = 0.000040369243
$ws.Cells.Item(34,4).Value = 1.352902466666667
$ws.Cells.Item(34,5).Value = 52.370959
$ws.Cells.Item(34,6).Value = 0.000037722511
$ws.Cells.Item(35,4).Value = 0.8410644666666667
$ws.Cells.Item(35,5).Value = 52.899788
$ws.Cells.Item(35,6).Value = 0.00002345111
$ws.Cells.Item(36,4).Value = 1.040556
$ws.Cells.Item(36,5).Value = 58.14514
$ws.Cells.Item(36,6).Value = 0.000029013462
$ws.Cells.Item(37,4).Value = 0.5245417333333334
$ws.Cells.Item(37,5).Value = 51.893882
$ws.Cells.Item(37,6).Value = 0.000014625615
$ws.Cells.Item(38,4).Value = 0.3777955066666667
$ws.Cells.Item(38,5).Value = 49.66533
$ws.Cells.Item(38,6).Value = 0.000010533941
$ws.Cells.Item(39,4).Value = 0.07512338666666667
$ws.Cells.Item(39,5).Value = 41.46553
$ws.Cells.Item(39,6).Value = 0.0000020946394
$ws.Cells.Item(40,4).Value = 3.569020666666667
$ws.Cells.Item(40,5).Value = 41.52611
$ws.Cells.Item(40,6).Value = 0.000099513766
$ws.Cells.Item(41,4).Value = 0.1698388933333333
$ws.Cells.Item(41,5).Value = 35.169491
$ws.Cells.Item(41,6).Value = 0.000004735559
$ws.Cells.Item(42,4).Value = 2.749582066666667
$ws.Cells.Item(42,5).Value = 35.320942
$ws.Cells.Item(42,6).Value = 0.000076665644
$ws.Cells.Item(43,4).Value = 0.3140000000000001
$ws.Cells.Item(43,5).Value = 0
$ws.Cells.Item(43,6).Value = 0.0000087551532
$ws.Cells.Item(44,4).Value = 0.045172414
$ws.Cells.Item(44,5).Value = 10.071461
$ws.Cells.Item(44,6).Value = 0.0000012595268
$ws.Cells.Item(45,4).Value = 0.02260702533333334
$ws.Cells.Item(45,5).Value = 0.27844853
$ws.Cells.Item(45,6).Value = 0.00000063034386
$ws.Cells.Item(46,4).Value = 0.35
$ws.Cells.Item(46,5).Value = 52.756089
$ws.Cells.Item(46,6).Value = 0.0000097589287
$ws.Cells.Item(47,4).Value = 0.35
$ws.Cells.Item(47,5).Value = 52.756089
$ws.Cells.Item(47,6).Value = 0.0000097589287
$ws.Cells.Item(48,4).Value = 3.36
$ws.Cells.Item(48,5).Value = 52.756089
$ws.Cells.Item(48,6).Value = 0.000093685715
$ws.Cells.Item(49,4).Value = 1.1607132
$ws.Cells.Item(49,5).Value = 57.94436
$ws.Cells.Item(49,6).Value = 0.000032363764
$ws.Cells.Item(50,4).Value = 3.971818133333334
$ws.Cells.Item(50,5).Value = 58.84974
$ws.Cells.Item(50,6).Value = 0.00011074483
$ws.Cells.Item(51,4).Value = 1.1607132
$ws.Cells.Item(51,5).Value = 57.94436
$ws.Cells.Item(51,6).Value = 0.000032363764
$ws.Cells.Item(52,4).Value = 3.971818133333334
$ws.Cells.Item(52,5).Value = 58.84974
$ws.Cells.Item(52,6).Value = 0.00011074483
$ws.Cells.Item(53,4).Value = 1.021207333333333
$ws.Cells.Item(53,5).Value = 34.033436
$ws.Cells.Item(53,6).Value = 0.00002847397
$ws.Cells.Item(54,4).Value = 3.016965866666667
$ws.Cells.Item(54,5).Value = 27.226749
$ws.Cells.Item(54,6).Value = 0.000084121014
$ws.Cells.Item(55,4).Value = 3.1285706
$ws.Cells.Item(55,5).Value = 46.355488
$ws.Cells.Item(55,6).Value = 0.000087232849
$ws.Cells.Item(56,4).Value = 2.8738048
$ws.Cells.Item(56,5).Value = 42.668768
$ws.Cells.Item(56,6).Value = 0.000080129303
$ws.Cells.Item(57,4).Value = 0.19
$ws.Cells.Item(57,5).Value = 52.745617
$ws.Cells.Item(57,6).Value = 0.0000052977041
$ws.Cells.Item(58,4).Value = 0.274
$ws.Cells.Item(58,5).Value = 51.000662
$ws.Cells.Item(58,6).Value = 0.000007639847
$ws.Cells.Item(59,4).Value = 0.482
$ws.Cells.Item(59,5).Value = 56.301312
$ws.Cells.Item(59,6).Value = 0.000013439439
$ws.Cells.Item(60,4).Value = 3.622
$ws.Cells.Item(60,5).Value = 56.301312
$ws.Cells.Item(60,6).Value = 0.00010099097
$ws.Cells.Item(61,4).Value = 0.08509560000000001
$ws.Cells.Item(61,5).Value = 1.3240718
$ws.Cells.Item(61,6).Value = 0.0000023726911
$ws.Cells.Item(62,4).Value = 3.595603466666667
$ws.Cells.Item(62,5).Value = 47.209177
$ws.Cells.Item(62,6).Value = 0.00010025496
$ws.Cells.Item(63,4).Value = 3.595603466666667
$ws.Cells.Item(63,5).Value = 47.209177
$ws.Cells.Item(63,6).Value = 0.00010025496
$ws.Cells.Item(64,4).Value = 0.398
$ws.Cells.Item(64,5).Value = 55.613552
$ws.Cells.Item(64,6).Value = 0.000011097296
$ws.Cells.Item(65,4).Value = 3.538
$ws.Cells.Item(65,5).Value = 55.613552
$ws.Cells.Item(65,6).Value = 0.000098648828
$ws.Cells.Item(66,4).Value = 3.955000000000001
$ws.Cells.Item(66,5).Value = 59.224905
$ws.Cells.Item(66,6).Value = 0.00011027589
$ws.Cells.Item(67,4).Value = 0.8150000000000001
$ws.Cells.Item(67,5).Value = 59.224905
$ws.Cells.Item(67,6).Value = 0.000022724363
$ws.Cells.Item(68,4).Value = 0.705
$ws.Cells.Item(68,5).Value = 58.164787
$ws.Cells.Item(68,6).Value = 0.000019657271
$ws.Cells.Item(69,4).Value = 2839.49
$ws.Cells.Item(69,5).Value = 3661126.4
$ws.Cells.Item(69,6).Value = 0.079172515

# --- Add header column-description comments ---
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
